$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (recordUuid), shifting recordUuid -> F, userEmail -> G
$ws.Columns("E").Insert()

# New header for column E, matching the style of the other header cells
$ws.Range("E1").Value = "comment"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# New sample data for column E, rows 2-4
$ws.Range("E2").Value = "Sample Text"
$ws.Range("E3").Value = "Sample Text"
$ws.Range("E4").Value = "Sample Text"

# Update effectiveDate values in column B
$ws.Range("B2").Value = "2026-01-05 00:22:38"
$ws.Range("B3").Value = "2026-01-05 00:22:38"
$ws.Range("B4").Value = "2026-01-05 00:22:38"
